$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 17-18, pushing the "Fix on Explorer" item
# (and everything below it) down by two rows, preserving the existing
# one-blank-row-between-items layout.
$ws.Range("B17:C18").EntireRow.Insert()

# Fill the newly-freed row 17 with the new ToDo item.
$ws.Range("B17").Value = "Change Highstocks with Highcharts"
$ws.Range("C17").Value = "Open"

# Match the new selection recorded in the workbook.
$ws.Range("B18").Select()
